# Scenario_Component_Behavior.xlsx — settings tweak so the model is feasible
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OperationScenario_Behavior")

# hot_water_demand_annual (column I): 2,500,000 -> 2,700
$ws.Range("I2").Value = 2700

# appliance_electricity_demand_annual (column L): 4,000,000 -> 1,000
$ws.Range("L2").Value = 1000

# Reflect the saved selection/scroll state from the authored workbook
$ws.Activate()
$ws.Range("I3").Select()
